$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts old C..F to D..G)
$ws.Columns("C").Insert()

# Header for the new column
$ws.Range("C1").Value = "P"

# Fill the new column with the pressure value for each data row
$ws.Range("C2:C11").Value = 101325

$ws.Range("C14").Select()
